# A.schedule.xlsx - "add lc.42, more understanding required"
#
# Sheet2!B6 and Sheet2!B7 hold the raw lottery-number lists that Sheet1
# pulls in via =Sheet2!Bn formulas (and the historical "last 4 rows"
# lookups). The edit folds the old "20,42,84" entry into "42" (its own
# row) and "20,84,25,239" (the following row), so:
#   Sheet2!B6 : 25,239     -> 42
#   Sheet2!B7 : 20,42,84   -> 20,84,25,239
# Sheet1's formula cells referencing these (B6/B7, C6/C7, C8, D14/D15)
# recalc automatically. Two Sheet1 cells also pick up the highlighted
# ("duplicate found") fill that's already used on B2:B5 (C5 and B6).

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# --- update the underlying data on Sheet2 ---------------------------------
$ws2.Range("B6").Value = "42"
$ws2.Range("B7").Value = "20,84,25,239"

# --- apply the highlighted fill (same format as B2) to C5 and B6 ----------
$ws1.Range("B2").Copy()
$ws1.Range("C5").PasteSpecial(-4122)
$ws1.Range("B2").Copy()
$ws1.Range("B6").PasteSpecial(-4122)

# --- update selections / view state ----------------------------------------
# Sheet2's selection moves to B8 (and gains a 120% zoom like Sheet1 already
# has); Sheet1 is reselected last so it stays the active/visible tab with
# its own selection moved to F16.
$ws2.Activate()
$ws2.Range("B8").Select()
$excel.ActiveWindow.Zoom = 120

$ws1.Activate()
$ws1.Range("F16").Select()
